$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the spelling "Visualisations" -> "Visualizations" in the
#    "More Visualisations" bullet under "Before Meeting".
#    (Scope the replace to that single paragraph so the other
#    "Visualisations" occurrence earlier in the document is untouched.)
# ------------------------------------------------------------------
$moreVisPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "More Visualisations*") {
        $moreVisPara = $p
    }
}

$moreVisStart = $moreVisPara.Range.Start
$visWordRange = $d.Range($moreVisStart + 5, $moreVisStart + 19)
$visWordRange.Text = "Visualizations"

# ------------------------------------------------------------------
# 2) Locate the "Talk Plan" / "Generate Images" / "Learning
#    Slurm/Baskerville" bullets that follow, delete the block (which
#    removes the old "Generate Images" bullet entirely) and retype the
#    replacement bullets in the new order:
#       Talk Plan
#       Learning Slurm/Baskerville
#       Train 2 Models, one small one for showing the VAE generates clear Images
#       Train a second on the full dataset overnight to show it can learn
#       Generate Example Images
# ------------------------------------------------------------------
$talkPlanPara = $null
$learningPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Talk Plan*") {
        $talkPlanPara = $p
    }
    if ($p.Range.Text -like "Learning Slurm/Baskerville*") {
        $learningPara = $p
    }
}

$blockRange = $d.Range($talkPlanPara.Range.Start, $learningPara.Range.End)
$blockRange.Delete()

$newBullets = @(
    "Talk Plan",
    "Learning Slurm/Baskerville",
    "Train 2 Models, one small one for showing the VAE generates clear Images",
    "Train a second on the full dataset overnight to show it can learn",
    "Generate Example Images"
)

foreach ($bullet in $newBullets) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lr = $lastPara.Range
    $lr.Collapse(0)
    $lr.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.InsertAfter($bullet)
}
